# Updated symbol list on Sun Jan 22 05:47:18 UTC 2023 with GitHub Actions
# Refresh Price (column D) and Volume/1h-change (column E) figures for the
# crypto ticker rows. Values are assigned with a leading apostrophe so
# Excel stores them as literal text (matching the sheet's existing
# inlineStr/text cells) instead of re-interpreting them as numbers or
# percentages, which would silently change precision/formatting
# (e.g. trailing zeros, scientific notation).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'301.48"
$ws.Range("E2").Value = "'-0.65%"
$ws.Range("D3").Value = "'37.75"
$ws.Range("E3").Value = "'8.40%"
$ws.Range("D4").Value = "'5.012"
$ws.Range("E4").Value = "'-2.83%"
$ws.Range("D5").Value = "'0.07872"
$ws.Range("E5").Value = "'1.54%"
$ws.Range("D6").Value = "'2.193"
$ws.Range("E6").Value = "'-7.25%"
$ws.Range("D7").Value = "'8.014"
$ws.Range("E7").Value = "'-0.02%"
$ws.Range("D8").Value = "'4.002"
$ws.Range("E8").Value = "'1.36%"
$ws.Range("D9").Value = "'0.9107"
$ws.Range("E9").Value = "'-1.97%"
$ws.Range("D10").Value = "'0.1873"
$ws.Range("E10").Value = "'3.91%"
$ws.Range("D11").Value = "'0.09225"
$ws.Range("E11").Value = "'-7.04%"
$ws.Range("D12").Value = "'0.08428"
$ws.Range("E12").Value = "'-2.56%"
$ws.Range("D13").Value = "'0.03514"
$ws.Range("E13").Value = "'5.95%"
$ws.Range("D14").Value = "'0.09936"
$ws.Range("E14").Value = "'0.49%"
$ws.Range("D15").Value = "'0.001471"
$ws.Range("E15").Value = "'-1.61%"
$ws.Range("D16").Value = "'0.005624"
$ws.Range("E16").Value = "'-2.30%"
$ws.Range("D17").Value = "'3.477"
$ws.Range("E17").Value = "'0.30%"
$ws.Range("E18").Value = "'-1.91%"
$ws.Range("D20").Value = "'0.1298"
$ws.Range("E20").Value = "'-2.67%"
$ws.Range("D21").Value = "'4.566"
$ws.Range("E21").Value = "'4.84%"
$ws.Range("E22").Value = "'-2.61%"
$ws.Range("D23").Value = "'0.04643"
$ws.Range("E23").Value = "'1.44%"
$ws.Range("D24").Value = "'0.001228"
$ws.Range("E24").Value = "'0.84%"
$ws.Range("D25").Value = "'0.004450"
$ws.Range("E25").Value = "'-0.28%"
$ws.Range("D26").Value = "'0.0001299"
$ws.Range("E26").Value = "'-0.11%"
$ws.Range("D27").Value = "'0.0004742"
$ws.Range("E27").Value = "'39.75%"
$ws.Range("D39").Value = "'0.01750"
$ws.Range("E39").Value = "'-2.14%"
$ws.Range("D40").Value = "'0.04715"
$ws.Range("E40").Value = "'-1.70%"
$ws.Range("D41").Value = "'0.007860"
$ws.Range("E41").Value = "'1.38%"
$ws.Range("E42").Value = "'-1.40%"
$ws.Range("D43").Value = "'0.007650"
$ws.Range("E43").Value = "'6.77%"
$ws.Range("D44").Value = "'0.002288"
$ws.Range("E44").Value = "'8.99%"
$ws.Range("D45").Value = "'0.01087"
$ws.Range("E45").Value = "'18.41%"
$ws.Range("D46").Value = "'0.00006077"
$ws.Range("E46").Value = "'-0.75%"
$ws.Range("D47").Value = "'0.00000000749"
$ws.Range("E47").Value = "'-0.19%"
$ws.Range("D48").Value = "'8.673"
$ws.Range("E48").Value = "'183.18%"
$ws.Range("E49").Value = "'34.86%"
$ws.Range("D50").Value = "'0.00002096"
$ws.Range("E50").Value = "'-0.19%"
$ws.Range("D51").Value = "'0.0001997"
$ws.Range("E51").Value = "'-0.19%"
